# Rename the worksheet: "Usuario" -> "Datos Exportados"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Datos Exportados"

# The header row A1:B1 was a merged cell ("CARACTERISTICAS GENERALES").
# Unmerge it so each header cell can hold its own value/style.
$ws.Range("A1:B1").UnMerge()

# Column widths: A=26, B=17, C=17, D=25 (stored OOXML width, not the
# padded "ColumnWidth" the COM layer otherwise reports - subtract the
# standard 11/12 character padding so the persisted width matches).
$ws.Columns.Item(1).ColumnWidth = 26 - 11 / 12
$ws.Columns.Item(2).ColumnWidth = 17 - 11 / 12
$ws.Columns.Item(3).ColumnWidth = 17 - 11 / 12
$ws.Columns.Item(4).ColumnWidth = 25 - 11 / 12

# Extend the existing alternating row styles (white / light-grey fill,
# centered) from columns A:B into the new C:D columns, row by row, by
# copying the already-styled A:B cells across. This reuses the workbook's
# existing style indices instead of fabricating new ones.
$ws.Range("A2:B2").Copy()
$ws.Range("C2:D2").PasteSpecial(-4122)
$ws.Range("A3:B3").Copy()
$ws.Range("C3:D3").PasteSpecial(-4122)
$ws.Range("A4:B4").Copy()
$ws.Range("C4:D4").PasteSpecial(-4122)
$ws.Range("A5:B5").Copy()
$ws.Range("C5:D5").PasteSpecial(-4122)
$ws.Range("A6:B6").Copy()
$ws.Range("C6:D6").PasteSpecial(-4122)
$ws.Range("A7:B7").Copy()
$ws.Range("C7:D7").PasteSpecial(-4122)

# New row 8: copy the "white" row style (row 6, same as rows 2/4) down.
$ws.Range("A6:D6").Copy()
$ws.Range("A8:D8").PasteSpecial(-4122)

# Extend the header style (A1, already bold/green/white after unmerge)
# across C1:D1.
$ws.Range("A1").Copy()
$ws.Range("C1:D1").PasteSpecial(-4122)

# --- Header row ---
$ws.Range("A1").Value = "NOMBRE"
$ws.Range("B1").Value = "STOCK ACTUAL"
$ws.Range("C1").Value = "MÁXIMO STOCK"
$ws.Range("D1").Value = "ÚLTIMA ACTUALIZACIÓN"

# --- Data rows ---
$ws.Range("A2").Value = "asdfsdfs"
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = "22/11/2024"

$ws.Range("A3").Value = "asdfsdfs"
$ws.Range("B3").Value = 5
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = "22/11/2024"

$ws.Range("A4").Value = "asdfsdfs"
$ws.Range("B4").Value = 155
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = "22/11/2024"

$ws.Range("A5").Value = "Inventario Principal"
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 500
$ws.Range("D5").Value = "22/11/2024"

$ws.Range("A6").Value = "Inventario Secundario"
$ws.Range("B6").Value = 535
$ws.Range("C6").Value = 600
$ws.Range("D6").Value = "27/11/2024"

$ws.Range("A7").Value = "Patio"
$ws.Range("B7").Value = 355
$ws.Range("C7").Value = 1400
$ws.Range("D7").Value = "27/11/2024"

$ws.Range("A8").Value = "asdd"
$ws.Range("B8").Value = 406
$ws.Range("C8").Value = 1500
$ws.Range("D8").Value = "27/11/2024"

Write-Output "done"
